$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39 - this shifts existing rows 39..111 down to 40..112
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new observation
$ws.Range("A39").Value = 11
$ws.Range("B39").Value = "Vega Monumental Concepción"
$ws.Range("C39").Value = "Bíobío"
$ws.Range("D39").Value = 44883
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 100112024
$ws.Range("G39").Value = "Choclo"
$ws.Range("H39").Value = "Choclero"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 25000
$ws.Range("L39").Value = 26000
$ws.Range("M39").Value = 25500
$ws.Range("N39").Value = "$/malla 50 unidades"
$ws.Range("O39").Value = "Región de Arica y Parinacota"
$ws.Range("P39").Value = 510
$ws.Range("Q39").Value = 50
$ws.Range("R39").Value = "Hortaliza"
